$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.116.72'
$ws.Range("E2").Value = '  +3.43%  '
$ws.Range("D3").Value = '3.199.12'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'538.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").Value = "'144.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.00%  '
$ws.Range("D9").Value = "'7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = '  +3.89%  '
$ws.Range("D11").Value = "'0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '3.753.76'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").Value = "'26.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = '60.163.57'
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").Value = '3.216.47'
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("D18").Value = "'6.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = "'13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").Value = "'383.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D24").Value = "'70.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").Value = "'8.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.51%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '0.0₃0905'
$ws.Range("E28").Value = '  +2.63%  '
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = "'6.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = "'22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.87%  '
$ws.Range("D32").Value = "'5.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.74%  '
$ws.Range("E33").Value = '  +5.12%  '
$ws.Range("D34").Value = "'6.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.28%  '
$ws.Range("D35").Value = "'156.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.45%  '
$ws.Range("D36").Value = "'1.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").Value = '2.794.13'
$ws.Range("E37").Value = '  +5.66%  '
$ws.Range("D38").Value = "'25.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("D39").Value = "'0.0710"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.55%  '
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("D41").Value = "'4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = "'39.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("D43").Value = "'0.727"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("E44").Value = '  +4.89%  '
$ws.Range("D45").Value = '3.243.03'
$ws.Range("E45").Value = '  +1.94%  '
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").Value = "'6.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").Value = "'0.797"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.94%  '
$ws.Range("D50").Value = "'20.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
